# Adapt column header formatting to respective input file names (FV2404 / FV2410),
# turn the data range into a proper Excel Table ("Table1"), and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the 21 header cells in row 1 -------------------------------
# Columns A-J  : "<name>_old" -> "<name>_FV2404"
# Column  K    : "diff"       -> unchanged
# Columns L-U  : "<name>_new" -> "<name>_FV2410"
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn A1:U85 into an Excel Table with autofilter --------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U85"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (row 1) --------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
